$d = $word.ActiveDocument

# 1) Merge the three runs around "chung,...)." into a single run and drop the
#    proofErr (grammar) markers that bracketed "chung,...". A plain
#    Find/Replace over that exact span collapses the identically-formatted
#    runs into one run and removes the now-orphaned proofErr elements.
$rng = $d.Content
$found = $rng.Find.Execute(
    "Mục đích của tính toán này là hạn chế biến dạng của nền, móng và kết cấu trên móng trong phạm vi đảm bảo không xảy ra tình hình cản trở việc sử dụng bình thường của nhà và công trình nói chung, hay của từng kết cấu hoặc giảm tính bền vững lâu dài của chúng do xuất hiện các chuyển vị không cho phép (độ lún, nghiêng, thay đổi cao độ thiết kế và vị trí kết cấu, phá hoại các liên kết của chúng,…).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Mục đích của tính toán này là hạn chế biến dạng của nền, móng và kết cấu trên móng trong phạm vi đảm bảo không xảy ra tình hình cản trở việc sử dụng bình thường của nhà và công trình nói chung, hay của từng kết cấu hoặc giảm tính bền vững lâu dài của chúng do xuất hiện các chuyển vị không cho phép (độ lún, nghiêng, thay đổi cao độ thiết kế và vị trí kết cấu, phá hoại các liên kết của chúng,…).",
    2)

# 2) Remove both comments. Capture the location of the second comment
#    (around "A, B, D") before deleting it so we can drop a `_GoBack`
#    bookmark exactly where its closing mark used to be.
$goBackPos = $null
if ($d.Comments.Count -ge 2) {
    $goBackPos = $d.Comments.Item(2).Scope.End
}

while ($d.Comments.Count -gt 0) {
    $d.Comments.Item(1).Delete()
}

# 3) Re-anchor the `_GoBack` bookmark at the former comment-2 location.
#    Word only ever keeps a single `_GoBack` bookmark, so adding a new one
#    moves it (and removes the stale one that otherwise sits at the end of
#    the document).
if ($goBackPos -ne $null) {
    $goBackRange = $d.Range($goBackPos, $goBackPos)
    $d.Bookmarks.Add("_GoBack", $goBackRange)
}
